# Apply the crypto price/volume updates captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-detects a plain numeric-looking string (e.g. "559.71") as a
# Number when assigned through .Value. Prefixing it with an apostrophe -
# exactly what a person typing it into the grid would do to force text -
# keeps the cell stored as Text, matching the original inline-string cells
# that merely "look like" numbers (prices such as "64.334.69" use '.' as a
# thousands separator and are never ambiguous, so they need no prefix).
function Set-TextValue($range, [string]$text) {
    if ($text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

$updates = @(
    @{ Row = 2; D = "64.384.31"; E = "  +1.39%  " }
    @{ Row = 3; D = "3.094.19"; E = "  +0.87%  " }
    @{ Row = 4; D = $null; E = "  -0.01%  " }
    @{ Row = 5; D = "559.71"; E = "  +1.70%  " }
    @{ Row = 6; D = "144.69"; E = "  +3.18%  " }
    @{ Row = 7; D = $null; E = "  -0.03%  " }
    @{ Row = 8; D = "3.091.82"; E = "  +0.94%  " }
    @{ Row = 9; D = "0.506"; E = "  +0.79%  " }
    @{ Row = 10; D = $null; E = "  +0.88%  " }
    @{ Row = 11; D = "6.14"; E = "  -5.95%  " }
    @{ Row = 12; D = $null; E = "  +3.43%  " }
    @{ Row = 13; D = $null; E = "  +0.09%  " }
    @{ Row = 14; D = "35.21"; E = "  +0.75%  " }
    @{ Row = 15; D = "3.591.47"; E = "  +0.75%  " }
    @{ Row = 16; D = "64.356.81"; E = "  +1.39%  " }
    @{ Row = 17; D = "3.091.04"; E = "  +0.84%  " }
    @{ Row = 18; D = $null; E = "  +1.19%  " }
    @{ Row = 19; D = $null; E = "  -0.35%  " }
    @{ Row = 20; D = "485.90"; E = "  -0.07%  " }
    @{ Row = 21; D = "14.01"; E = "  +1.57%  " }
    @{ Row = 22; D = $null; E = "  -0.24%  " }
    @{ Row = 23; D = "7.56"; E = "  +4.04%  " }
    @{ Row = 24; D = "14.26"; E = "  +12.68%  " }
    @{ Row = 25; D = "81.37"; E = "  +0.37%  " }
    @{ Row = 26; D = $null; E = "  +0.15%  " }
    @{ Row = 27; D = $null; E = "  +1.18%  " }
    @{ Row = 28; D = "8.04"; E = "  +1.42%  " }
    @{ Row = 29; D = $null; E = "  +2.72%  " }
    @{ Row = 30; D = $null; E = "  -0.01%  " }
    @{ Row = 31; D = $null; E = "  +0.59%  " }
    @{ Row = 32; D = "1.15"; E = "  -1.20%  " }
    @{ Row = 33; D = $null; E = "  +1.20%  " }
    @{ Row = 34; D = $null; E = "  -1.68%  " }
    @{ Row = 35; D = "6.24"; E = "  +4.12%  " }
    @{ Row = 36; D = "55.74"; E = "  +0.26%  " }
    @{ Row = 37; D = "3.01"; E = "  +16.97%  " }
    @{ Row = 38; D = "454.83"; E = "  -2.45%  " }
    @{ Row = 39; D = "0.0411"; E = "  +3.06%  " }
    @{ Row = 40; D = "0.0818"; E = "  -0.71%  " }
    @{ Row = 41; D = "2.969.21"; E = "  -2.91%  " }
    @{ Row = 42; D = $null; E = "  -0.27%  " }
    @{ Row = 43; D = $null; E = "  -5.66%  " }
    @{ Row = 44; D = "28.28"; E = "  +0.43%  " }
    @{ Row = 45; D = $null; E = "  +2.01%  " }
    @{ Row = 46; D = $null; E = "  -0.03%  " }
    @{ Row = 47; D = $null; E = "  +3.97%  " }
    @{ Row = 48; D = $null; E = "  +1.84%  " }
    @{ Row = 49; D = "118.54"; E = "  +1.22%  " }
    @{ Row = 50; D = "0.0₃0518"; E = "  +1.07%  " }
    @{ Row = 51; D = $null; E = "  +0.34%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) { Set-TextValue $ws.Range("D$($u.Row)") $u.D }
    if ($null -ne $u.E) { Set-TextValue $ws.Range("E$($u.Row)") $u.E }
}
